# Slide 9 ("random-walk.9") contains a probability-tree diagram built from
# grouped Oval shapes.  Two of the ovals (inside the group named "Group 1")
# show the running sequence of coin-toss outcomes:
#   - Oval id=32 currently reads "HH" (solid green)
#   - Oval id=42 currently reads "H" (green) + "T" (red)
#
# The edit recolors/retexts them so that:
#   - Oval id=32 becomes "T" (red) followed by "H" (green)
#   - Oval id=42 becomes a single "TT" run (red)

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(9)
$grp = $s.Shapes.Item(7)          # "Group 1"

# --- Oval 35 (id=32): "HH" (green) -> "T" (red) + "H" (green) ---
$ovalHH = $grp.GroupItems.Item(5)
$trHH = $ovalHH.TextFrame.TextRange

$firstChar = $trHH.Characters(1,1)
$firstChar.Font.Color.RGB = 255      # COLORREF 0x0000FF -> srgbClr FF0000 (red)
$firstChar.Text = "T"

# --- Oval 35 (id=42): "H" (green) + "T" (red) -> single "TT" run (red) ---
$ovalHT = $grp.GroupItems.Item(7)
$trHT = $ovalHT.TextFrame.TextRange

# Delete the leading "H" character; this collapses the text range back down
# to the single remaining (red) run instead of leaving two runs behind.
$lead = $trHT.Characters(1,1)
$lead.Text = ""

# Now there is one run left ("T", red). Growing it back to "TT" keeps it as
# a single run with the same (red) formatting.
$trHT.Text = "TT"
